$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.420.59"
$ws.Range("E2").Value = "  +4.06%  "

$ws.Range("D3").Value = "2.458.26"
$ws.Range("E3").Value = "  +1.50%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "321.56"
$ws.Range("E5").Value = "  +2.16%  "

$ws.Range("D6").Value = "105.68"
$ws.Range("E6").Value = "  +4.55%  "

$ws.Range("E7").Value = "  +1.19%  "

$ws.Range("D9").Value = "0.538"
$ws.Range("E9").Value = "  +3.04%  "

$ws.Range("D10").Value = "36.20"
$ws.Range("E10").Value = "  +2.22%  "

$ws.Range("E12").Value = "  +0.57%  "

$ws.Range("D13").Value = "18.40"
$ws.Range("E13").Value = "  -3.02%  "

$ws.Range("D14").Value = "7.09"
$ws.Range("E14").Value = "  +1.71%  "

$ws.Range("D15").Value = "2.848.52"
$ws.Range("E15").Value = "  +1.69%  "

$ws.Range("D16").Value = "2.457.35"
$ws.Range("E16").Value = "  +1.37%  "

$ws.Range("D17").Value = "0.846"
$ws.Range("E17").Value = "  +1.69%  "

$ws.Range("D18").Value = "46.267.34"
$ws.Range("E18").Value = "  +4.05%  "

$ws.Range("D19").Value = "12.75"
$ws.Range("E19").Value = "  +2.25%  "

$ws.Range("E20").Value = "  +0.66%  "

$ws.Range("D21").Value = "0.0₃0937"
$ws.Range("E21").Value = "  +1.56%  "

$ws.Range("D22").Value = "70.61"
$ws.Range("E22").Value = "  +2.74%  "

$ws.Range("E23").Value = "  +4.63%  "

$ws.Range("D24").Value = "247.93"
$ws.Range("E24").Value = "  +2.25%  "

$ws.Range("E25").Value = "  +1.72%  "

$ws.Range("D26").Value = "26.16"
$ws.Range("E26").Value = "  +3.70%  "

$ws.Range("E27").Value = "  +0.04%  "

$ws.Range("D28").Value = "2.29"
$ws.Range("E28").Value = "  +0.24%  "

$ws.Range("D29").Value = "9.79"
$ws.Range("E29").Value = "  +2.61%  "

$ws.Range("D30").Value = "34.77"
$ws.Range("E30").Value = "  +4.53%  "

$ws.Range("D31").Value = "49.51"
$ws.Range("E31").Value = "  +2.16%  "

$ws.Range("E32").Value = "  +3.23%  "

$ws.Range("D33").Value = "19.82"
$ws.Range("E33").Value = "  +2.21%  "

$ws.Range("E34").Value = "  +3.24%  "

$ws.Range("D36").Value = "0.0765"
$ws.Range("E36").Value = "  -1.38%  "

$ws.Range("D37").Value = "4.60"
$ws.Range("E37").Value = "  +2.23%  "

$ws.Range("E38").Value = "  +0.76%  "

$ws.Range("E39").Value = "  +3.21%  "

$ws.Range("D40").Value = "122.82"
$ws.Range("E40").Value = "  +2.41%  "

$ws.Range("E41").Value = "  +1.85%  "

$ws.Range("E42").Value = "  +1.72%  "

$ws.Range("D43").Value = "20.88"
$ws.Range("E43").Value = "  -0.76%  "

$ws.Range("E44").Value = "  +0.86%  "

$ws.Range("D45").Value = "1.978.43"
$ws.Range("E45").Value = "  +1.94%  "

$ws.Range("D46").Value = "3.00"
$ws.Range("E46").Value = "  +1.81%  "

$ws.Range("D47").Value = "2.10"
$ws.Range("E47").Value = "  -3.33%  "

$ws.Range("D48").Value = "1.86"
$ws.Range("E48").Value = "  +11.33%  "

$ws.Range("D49").Value = "9.10"
$ws.Range("E49").Value = "  -3.67%  "

$ws.Range("D50").Value = "5.15"
$ws.Range("E50").Value = "  +10.86%  "

$ws.Range("D51").Value = "78.74"
$ws.Range("E51").Value = "  +4.15%  "
